# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Swap the full record (columns B..AC) between rows 168 & 169, and between rows 173 & 174.
# Column A (the sequential row index) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, [int]$row1, [int]$row2, [string]$firstCol, [string]$lastCol) {
    $rng1 = $ws.Range($firstCol + $row1 + ":" + $lastCol + $row1)
    $rng2 = $ws.Range($firstCol + $row2 + ":" + $lastCol + $row2)

    $cols = $rng1.Columns.Count

    for ($i = 1; $i -le $cols; $i++) {
        $c1 = $rng1.Cells.Item(1, $i)
        $c2 = $rng2.Cells.Item(1, $i)

        $v1 = $c1.Value2
        $v2 = $c2.Value2

        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

Swap-Rows $ws 168 169 "B" "AC"
Swap-Rows $ws 173 174 "B" "AC"
